$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.61%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'27.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'3.53%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.114"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.58%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05689"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'1.54%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.527"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.57%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8195"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'0.79%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8556"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'1.33%"
$ws.Range("E8").Style = "Normal"
$ws.Range("B9").Value = "'WazirX"
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").Value = "'0.1335"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.40%"
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'MandalaExchangeToken"
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'0.06942"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.78%"
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'BitrueCoin"
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'0.02862"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.94%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'BitMartToken"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.09385"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-0.12%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'BitForexToken"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.001521"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.06%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'CoinExToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.04075"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-12.61%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'One"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.0006031"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.53%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.006216"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.39%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'LEO"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.511"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.65%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'GateToken"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'3.010"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.26%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'BTSEToken"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'2.230"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'8.51%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'BitpandaEcosystemToken"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.3164"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.23%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.03208"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.34%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'-0.10%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.555"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-5.10%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'1.80%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001219"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-2.17%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004473"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-2.38%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.00009901"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'3.15%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'-25.21%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03724"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.64%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "'BKEXToken"
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = "'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").Value = "'0.1057"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-22.46%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "'CEJI"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'0.002440"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-8.25%"
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "'KickToken"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'0.003447"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-44.12%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009717"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'17.66%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005001"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-8.09%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'0.03%"
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'-8.15%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002510"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-3.00%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.03%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("E50").Style = "Normal"
